# Updates the cryptos price/volume table to the latest scraped snapshot.
# (commit: "Updated cryptos list on Wed Sep 11 07:50:00 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold plain-looking decimals (e.g. "132.30") that must stay
# TEXT (matching every other cell in the column) instead of being auto-coerced
# to a Number (which would silently drop the trailing zero). Force text via
# NumberFormat, write the value, then restore the default style so no other
# cell formatting changes.
$textCells = @("D5", "D6", "D11", "D13", "D19", "D21", "D23", "D24", "D28", "D29", "D31", "D32", "D37", "D38", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D51")
foreach ($cell in $textCells) { $ws.Range($cell).NumberFormat = "@" }

# Cell-by-cell updates (Coin / Link / Price / Volume(1h) columns).
$ws.Range('D2').Value = '56.494.71'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '2.332.76'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('D5').Value = '513.55'
$ws.Range('D6').Value = '132.30'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -1.06%  '
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('D11').Value = '5.29'
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '23.58'
$ws.Range('E13').Value = '  -0.91%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.745.81'
$ws.Range('E14').Value = '  -1.17%  '
$ws.Range('D15').Value = '56.492.66'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('D17').Value = '2.335.79'
$ws.Range('E17').Value = '  -1.61%  '
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '325.13'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('D21').Value = '6.70'
$ws.Range('E21').Value = '  +2.03%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '61.73'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').Value = '8.73'
$ws.Range('E24').Value = '  +11.20%  '
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('E27').Value = '  +4.29%  '
$ws.Range('D28').Value = '168.03'
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('D29').Value = '1.68'
$ws.Range('E29').Value = '  -1.41%  '
$ws.Range('D30').Value = '0.0₃0721'
$ws.Range('E30').Value = '  -3.59%  '
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('D32').Value = '18.38'
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').Value = '0.888'
$ws.Range('E37').Value = '  -4.83%  '
$ws.Range('D38').Value = '153.11'
$ws.Range('E38').Value = '  +11.52%  '
$ws.Range('E39').Value = '  +0.68%  '
$ws.Range('D40').Value = '38.48'
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').Value = '3.57'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').Value = '279.73'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('E44').Value = '  -2.02%  '
$ws.Range('D45').Value = '0.0929'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0495'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.558'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('D48').Value = '18.18'
$ws.Range('E48').Value = '  +4.77%  '
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').Value = '17.20'
$ws.Range('E51').Value = '  +1.74%  '

# Restore the original (default) style on the cells we forced to text above.
foreach ($cell in $textCells) { $ws.Range($cell).Style = "Normal" }

Write-Output "Updated 84 cells in the cryptos table"
